{"js": "// The document (\"ClassiTestate\" - \"tested classes\") lists classes that\n// have been unit-tested, one per bulleted paragraph (\"Paragrafoelenco\"\n// list style, numId 1): \"PawnClass\", then \"Position\".\n//\n// The commit adds two more finished classes to the bulleted list, in\n// this order, right after \"Position\":\n//   - \"PawnMap\"    (flagged by the proofer as a possible misspelling,\n//                    hence the surrounding <w:proofErr> spell-check\n//                    markers in the canonical XML)\n//   - \"Heuristic\"\n//\n// The trailing \"_GoBack\" bookmark (automatically maintained by Word to\n// remember the last edit position) moves from the end of the old last\n// paragraph (\"Position\") to the end of the new last paragraph\n// (\"Heuristic\").\n\nconst doc = context.document;\nconst body = doc.body;\n\n// 1. Drop the existing \"_GoBack\" bookmark; it currently sits at the end\n//    of the \"Position\" paragraph and needs to move to the new last\n//    paragraph once the new content is appended.\ndoc.deleteBookmark(\"_GoBack\");\n\n// 2. Find the current last paragraph (\"Position\") so the new paragraphs\n//    can be inserted right after it, inheriting its list formatting.\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// 3. Insert the two new bulleted paragraphs as raw OOXML so the\n//    \"PawnMap\" spell-check markers and the relocated bookmark can be\n//    reproduced exactly; both paragraphs reuse the same paragraph style,\n//    numbering (list) properties and run formatting (40 half-points /\n//    20pt) as the existing list items.\nconst insertionRange = lastParagraph.getRange(\"After\");\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragrafoelenco\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n              <w:t>PawnMap</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragrafoelenco\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n              <w:t>Heuristic</w:t>\n            </w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ninsertionRange.insertOoxml(ooxml, Word.InsertLocation.after);\nawait context.sync();\n", "ps1": "# The document (\"ClassiTestate\" - \"tested classes\") lists classes that\n# have been unit-tested, one per bulleted paragraph (\"Paragrafoelenco\"\n# list style, numId 1): \"PawnClass\", then \"Position\".\n#\n# This commit adds two more finished classes to the bulleted list, in\n# this order, right after \"Position\":\n#   - \"PawnMap\"    (flagged by the proofer as a possible misspelling,\n#                    hence the surrounding <w:proofErr> spell-check\n#                    markers in the canonical XML)\n#   - \"Heuristic\"\n#\n# The trailing \"_GoBack\" bookmark (automatically maintained by Word to\n# remember the last edit position) moves from the end of the old last\n# paragraph (\"Position\") to the end of the new last paragraph\n# (\"Heuristic\").\n\n$d = $word.ActiveDocument\n\n# 1. Drop the existing \"_GoBack\" bookmark; it currently sits at the end\n#    of the \"Position\" paragraph and needs to move to the new last\n#    paragraph once the new content is appended.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2. Collapse a range to the very end of the document content so the new\n#    paragraphs land after \"Position\" (the current last paragraph).\n$rng = $d.Content\n$rng.Collapse(0)\n\n# 3. Insert the two new bulleted paragraphs as raw WordprocessingML so the\n#    \"PawnMap\" spell-check markers and the relocated bookmark can be\n#    reproduced exactly; both paragraphs reuse the same paragraph style,\n#    numbering (list) properties and run formatting (40 half-points /\n#    20pt) as the existing list items.\n$xml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragrafoelenco\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:proofErr w:type=\"spellStart\"/>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n              <w:t>PawnMap</w:t>\n            </w:r>\n            <w:proofErr w:type=\"spellEnd\"/>\n          </w:p>\n          <w:p>\n            <w:pPr>\n              <w:pStyle w:val=\"Paragrafoelenco\"/>\n              <w:numPr>\n                <w:ilvl w:val=\"0\"/>\n                <w:numId w:val=\"1\"/>\n              </w:numPr>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n            </w:pPr>\n            <w:r>\n              <w:rPr>\n                <w:sz w:val=\"40\"/>\n              </w:rPr>\n              <w:t>Heuristic</w:t>\n            </w:r>\n            <w:bookmarkStart w:id=\"0\" w:name=\"_GoBack\"/>\n            <w:bookmarkEnd w:id=\"0\"/>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n$rng.InsertXML($xml)\n"}
